$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.612.14"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "2.219.38"
$ws.Range("E3").Value = "  +2.89%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'230.51"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'60.77"
$ws.Range("E7").Value = "  -3.18%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.401"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").Value = "'58.91"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  +5.60%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "2.547.21"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "'15.69"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'21.77"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "'5.56"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "2.229.22"
$ws.Range("E18").Value = "  +3.10%  "
$ws.Range("D19").Value = "41.498.85"
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("D20").Value = "'72.83"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D21").Value = "0.0₃0893"
$ws.Range("E21").Value = "  +5.05%  "
$ws.Range("D22").Value = "'6.03"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'250.92"
$ws.Range("E23").Value = "  +10.30%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "'2.34"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").Value = "'9.56"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").Value = "'167.30"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").Value = "'19.94"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "'4.95"
$ws.Range("E34").Value = "  +5.34%  "
$ws.Range("D35").Value = "'4.61"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").Value = "'0.0623"
$ws.Range("D37").Value = "'6.61"
$ws.Range("E37").Value = "  -4.99%  "
$ws.Range("D38").Value = "'3.66"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("D39").Value = "'2.36"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "'0.000250"
$ws.Range("E40").Value = "  +31.38%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0238"
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.79"
$ws.Range("E43").Value = "  -5.02%  "
$ws.Range("D44").Value = "'8.67"
$ws.Range("E44").Value = "  +10.54%  "
$ws.Range("D45").Value = "'0.0975"
$ws.Range("E45").Value = "  +5.98%  "
$ws.Range("D46").Value = "'98.47"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").Value = "1.467.89"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("E49").Value = "  -6.69%  "
$ws.Range("D50").Value = "'2.80"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'52.15"
$ws.Range("E51").Value = "  +7.19%  "
